# "Generate Report for Handoff"
#
# Updates the localization-status report:
#  - Status moves from "Handoff transform failed" to "Ready for handoff"
#    (Overview sheet + each language sheet).
#  - Each language sheet gets a new "Latest Handoff File" hyperlink (col C)
#    pointing at the freshly generated .xlf handoff package, plus a
#    "Latest Handoff Datetime" (col D) stamp.
#  - "Handoff Reason" (col H) flips from Ignored to Include for the file
#    that is now being handed off.

$wb = $excel.ActiveWorkbook

# Cornflowerblue (RGB 100,149,237 -> 0x6495ED) packed as BGR long, matching
# the workbook's existing HyperLink font/style used on column A.
$hyperlinkColor = 15570276

function Set-HandoffHyperlink($ws, $cellRef, $address, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText)
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
    $ws.Range($cellRef).Font.Underline = $true
}

# ---- Overview sheet: roll up the new status text ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Ready for handoff"
Set-HandoffHyperlink $zhcn "C2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/876c2a27fd5c10e8992efbd3de3fd177a6c91e8b/e2e/fc9bd0dc-745f-45f5-8096-bf9b461f5a24.876c2a27fd5c10e8992efbd3de3fd177a6c91e8b.zh-cn.xlf" `
    "fc9bd0dc-745f-45f5-8096-bf9b461f5a24.876c2a27fd5c10e8992efbd3de3fd177a6c91e8b.zh-cn.xlf"
$zhcn.Range("D2").Value = "2016-01-27 08:20:43"
$zhcn.Range("H2").Value = "Include"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Ready for handoff"
Set-HandoffHyperlink $dede "C2" `
    "https://github.com/OpenLocalizationTest/oltest/blob/876c2a27fd5c10e8992efbd3de3fd177a6c91e8b/e2e/fc9bd0dc-745f-45f5-8096-bf9b461f5a24.876c2a27fd5c10e8992efbd3de3fd177a6c91e8b.de-de.xlf" `
    "fc9bd0dc-745f-45f5-8096-bf9b461f5a24.876c2a27fd5c10e8992efbd3de3fd177a6c91e8b.de-de.xlf"
$dede.Range("D2").Value = "2016-01-27 08:20:58"
$dede.Range("H2").Value = "Include"
